$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Fix D8:D11 which were stored as text ("inlineStr") -> make them numeric ---
$ws.Cells.Item(8, 4).Value  = 500290
$ws.Cells.Item(9, 4).Value  = 500530
$ws.Cells.Item(10, 4).Value = 532500
$ws.Cells.Item(11, 4).Value = 532754

# --- Append new rows 12-14 with additional stock data ---

# Row 12: TATACHEM
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "TATACHEM"
$ws.Cells.Item(12, 3).Value = "Tata Chemicals Limited"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "500770"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = 2.34
$ws.Cells.Item(12, 6).Value = 1111.4
$ws.Cells.Item(12, 7).Value = 3432027
$ws.Cells.Item(12, 8).Value = "day"
$ws.Cells.Item(12, 9).Value = "11/06/2024 10:32:41"

# Row 13: HDFCLIFE
$ws.Cells.Item(13, 1).Value = 2
$ws.Cells.Item(13, 2).Value = "HDFCLIFE"
$ws.Cells.Item(13, 3).Value = "HDFC Life Insurance Company Ltd"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "540777"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = 0.44
$ws.Cells.Item(13, 6).Value = 571.7
$ws.Cells.Item(13, 7).Value = 6417546
$ws.Cells.Item(13, 8).Value = "day"
$ws.Cells.Item(13, 9).Value = "11/06/2024 10:32:41"

# Row 14: NATIONALUM
$ws.Cells.Item(14, 1).Value = 3
$ws.Cells.Item(14, 2).Value = "NATIONALUM"
$ws.Cells.Item(14, 3).Value = "National Aluminium Company Limited"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "532234"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = -0.13
$ws.Cells.Item(14, 6).Value = 183.17
$ws.Cells.Item(14, 7).Value = 16052453
$ws.Cells.Item(14, 8).Value = "day"
$ws.Cells.Item(14, 9).Value = "11/06/2024 10:32:41"
